# Update countries & provincias Spain
#
# 1. Swap the "Santa Lucia" / "Laos" pair and the "Islas Malvinas" /
#    "Groenlandia" pair (the underlying numbers for each pair are identical,
#    so the only visible effect is the two country labels trading places).
# 2. Bump the "last updated" timestamp from 10:12 to 11:29.
# 3. Refresh the daily case counters for a handful of countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Swap mislabeled rows -------------------------------------------
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Laos"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- 2. Update the "datos actualizados" timestamp -----------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 11:29"

# --- 3. Refresh the numeric counters ------------------------------------
# Row 20 - Banglades
$ws.Range("B20").Value = 137787
$ws.Range("C20").Value = 3809
$ws.Range("D20").Value = 55727
$ws.Range("E20").Value = 80322
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = 1738

# Row 32 - Indonesia
$ws.Range("B32").Value = 54010
$ws.Range("C32").Value = 1198
$ws.Range("D32").Value = 22936
$ws.Range("E32").Value = 28320
$ws.Range("G32").Value = 34
$ws.Range("H32").Value = 2754

# Row 40 - Oman
$ws.Range("B40").Value = 38150
$ws.Range("C40").Value = 1197
$ws.Range("D40").Value = 21200
$ws.Range("E40").Value = 16787
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 163

# Row 41 - Filipinas
$ws.Range("B41").Value = 35455
$ws.Range("C41").Value = 652
$ws.Range("D41").Value = 9686
$ws.Range("E41").Value = 24525
$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 1244

# Row 42 - Polonia (D unchanged)
$ws.Range("B42").Value = 33907
$ws.Range("C42").Value = 193
$ws.Range("E42").Value = 11921
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 1438

# Row 56 - Austria
$ws.Range("B56").Value = 17654
$ws.Range("C56").Value = 74
$ws.Range("D56").Value = 16401
$ws.Range("E56").Value = 551
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 702

# Row 60 - Moldavia (B, C unchanged)
$ws.Range("D60").Value = 9081
$ws.Range("E60").Value = 6473
$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 526

# Row 68 - Marruecos (F, G, H unchanged)
$ws.Range("B68").Value = 11986
$ws.Range("C68").Value = 109
$ws.Range("D68").Value = 8730
$ws.Range("E68").Value = 3036

# Row 73 - Malasia (F, G, H unchanged)
$ws.Range("B73").Value = 8634
$ws.Range("C73").Value = 18
$ws.Range("D73").Value = 8318
$ws.Range("E73").Value = 195

# Row 108 - Sri Lanka (only D, E changed)
$ws.Range("D108").Value = 1661
$ws.Range("E108").Value = 361

# Row 120 - Eslovenia
$ws.Range("B120").Value = 1581
$ws.Range("C120").Value = 9
$ws.Range("D120").Value = 1384
$ws.Range("E120").Value = 88

# Row 190 - Gambia (D unchanged)
$ws.Range("B190").Value = 45
$ws.Range("C190").Value = 1
$ws.Range("E190").Value = 17
